$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    3  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    4  = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 645.3272768299601, 676.4434635367506)
    5  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    6  = @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1.742940831014585)
    7  = @(0.3048080303191223, 0.3127903958511391, 0.8054896365839992, 8.660232485948974, 10.08332054870323)
    8  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    9  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    10 = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 8.660232485948974, 14.90378790461981)
    11 = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 0.496779210170732, 6.740334628841572)
    12 = @(0.6753301551942219, 1.667794583268128, 337.1190423067083, 645.3272768299601, 984.7894438751307)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
